$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "L1cam"
$ws.Cells.Item(2,3).Value = "Alcam"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 19.72083766666667
$ws.Cells.Item(2,8).Value = 59.162513
$ws.Cells.Item(2,9).Value = 0.8016210077351786
$ws.Cells.Item(2,10).Value = 0.8016210077351787
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 62.12558000000001
$ws.Cells.Item(2,14).Value = 186.37674
$ws.Cells.Item(2,15).Value = 0.9736910227596813
$ws.Cells.Item(2,16).Value = 0.9736910227596813
$ws.Cells.Item(2,17).Value = 1225.168478127513
$ws.Cells.Item(2,18).Value = 11026.51630314762
$ws.Cells.Item(2,19).Value = 0.7805311788873125
$ws.Cells.Item(2,20).Value = 0.7805311788873126

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "L1cam"
$ws.Cells.Item(3,3).Value = "Alcam"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 19.72083766666667
$ws.Cells.Item(3,8).Value = 59.162513
$ws.Cells.Item(3,9).Value = 0.8016210077351786
$ws.Cells.Item(3,10).Value = 0.8016210077351787
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.5683613333333334
$ws.Cells.Item(3,14).Value = 1.705084
$ws.Cells.Item(3,15).Value = 0.008907897969731461
$ws.Cells.Item(3,16).Value = 0.008907897969731461
$ws.Cells.Item(3,17).Value = 11.20856159067689
$ws.Cells.Item(3,18).Value = 100.877054316092
$ws.Cells.Item(3,19).Value = 0.007140758147298286
$ws.Cells.Item(3,20).Value = 0.007140758147298287

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "L1cam"
$ws.Cells.Item(4,3).Value = "Alcam"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 19.72083766666667
$ws.Cells.Item(4,8).Value = 59.162513
$ws.Cells.Item(4,9).Value = 0.8016210077351786
$ws.Cells.Item(4,10).Value = 0.8016210077351787
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.110262
$ws.Cells.Item(4,14).Value = 3.330786
$ws.Cells.Item(4,15).Value = 0.01740107927058724
$ws.Cells.Item(4,16).Value = 0.01740107927058724
$ws.Cells.Item(4,17).Value = 21.89529666946867
$ws.Cells.Item(4,18).Value = 197.057670025218
$ws.Cells.Item(4,19).Value = 0.01394907070056787
$ws.Cells.Item(4,20).Value = 0.01394907070056788

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "L1cam"
$ws.Cells.Item(5,3).Value = "Alcam"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.099159
$ws.Cells.Item(5,8).Value = 0.297477
$ws.Cells.Item(5,9).Value = 0.004030657259573097
$ws.Cells.Item(5,10).Value = 0.004030657259573097
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 62.12558000000001
$ws.Cells.Item(5,14).Value = 186.37674
$ws.Cells.Item(5,15).Value = 0.9736910227596813
$ws.Cells.Item(5,16).Value = 0.9736910227596813
$ws.Cells.Item(5,17).Value = 6.160310387220001
$ws.Cells.Item(5,18).Value = 55.44279348498
$ws.Cells.Item(5,19).Value = 0.003924614789467463
$ws.Cells.Item(5,20).Value = 0.003924614789467464

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "L1cam"
$ws.Cells.Item(6,3).Value = "Alcam"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.099159
$ws.Cells.Item(6,8).Value = 0.297477
$ws.Cells.Item(6,9).Value = 0.004030657259573097
$ws.Cells.Item(6,10).Value = 0.004030657259573097
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.5683613333333334
$ws.Cells.Item(6,14).Value = 1.705084
$ws.Cells.Item(6,15).Value = 0.008907897969731461
$ws.Cells.Item(6,16).Value = 0.008907897969731461
$ws.Cells.Item(6,17).Value = 0.056358141452
$ws.Cells.Item(6,18).Value = 0.507223273068
$ws.Cells.Item(6,19).Value = 0.00003590468361923456
$ws.Cells.Item(6,20).Value = 0.00003590468361923457

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "L1cam"
$ws.Cells.Item(7,3).Value = "Alcam"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.099159
$ws.Cells.Item(7,8).Value = 0.297477
$ws.Cells.Item(7,9).Value = 0.004030657259573097
$ws.Cells.Item(7,10).Value = 0.004030657259573097
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.110262
$ws.Cells.Item(7,14).Value = 3.330786
$ws.Cells.Item(7,15).Value = 0.01740107927058724
$ws.Cells.Item(7,16).Value = 0.01740107927058724
$ws.Cells.Item(7,17).Value = 0.110092469658
$ws.Cells.Item(7,18).Value = 0.990832226922
$ws.Cells.Item(7,19).Value = 0.00007013778648639939
$ws.Cells.Item(7,20).Value = 0.00007013778648639941

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "L1cam"
$ws.Cells.Item(8,3).Value = "Alcam"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 4.781202
$ws.Cells.Item(8,8).Value = 14.343606
$ws.Cells.Item(8,9).Value = 0.1943483350052483
$ws.Cells.Item(8,10).Value = 0.1943483350052483
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 62.12558000000001
$ws.Cells.Item(8,14).Value = 186.37674
$ws.Cells.Item(8,15).Value = 0.9736910227596813
$ws.Cells.Item(8,16).Value = 0.9736910227596813
$ws.Cells.Item(8,17).Value = 297.03494734716
$ws.Cells.Item(8,18).Value = 2673.31452612444
$ws.Cells.Item(8,19).Value = 0.1892352290829014
$ws.Cells.Item(8,20).Value = 0.1892352290829014

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "L1cam"
$ws.Cells.Item(9,3).Value = "Alcam"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 4.781202
$ws.Cells.Item(9,8).Value = 14.343606
$ws.Cells.Item(9,9).Value = 0.1943483350052483
$ws.Cells.Item(9,10).Value = 0.1943483350052483
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.5683613333333334
$ws.Cells.Item(9,14).Value = 1.705084
$ws.Cells.Item(9,15).Value = 0.008907897969731461
$ws.Cells.Item(9,16).Value = 0.008907897969731461
$ws.Cells.Item(9,17).Value = 2.717450343656
$ws.Cells.Item(9,18).Value = 24.457053092904
$ws.Cells.Item(9,19).Value = 0.001731235138813941
$ws.Cells.Item(9,20).Value = 0.001731235138813941

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "L1cam"
$ws.Cells.Item(10,3).Value = "Alcam"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 4.781202
$ws.Cells.Item(10,8).Value = 14.343606
$ws.Cells.Item(10,9).Value = 0.1943483350052483
$ws.Cells.Item(10,10).Value = 0.1943483350052483
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.110262
$ws.Cells.Item(10,14).Value = 3.330786
$ws.Cells.Item(10,15).Value = 0.01740107927058724
$ws.Cells.Item(10,16).Value = 0.01740107927058724
$ws.Cells.Item(10,17).Value = 5.308386894924001
$ws.Cells.Item(10,18).Value = 47.77548205431601
$ws.Cells.Item(10,19).Value = 0.00338187078353297
$ws.Cells.Item(10,20).Value = 0.003381870783532971

